# This script replays the daily "cryptos list" refresh performed by the
# GitHub Actions bot: it overwrites the Price (column D) and Volume(1h)
# (column E) figures for each coin row with freshly scraped values, and
# also fixes the order of the ARBITRUM / TrustWalletToken rows (40-41),
# which were swapped upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.497.06"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.43%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5090"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08345"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.81%  "
$ws.Range("E10").Value = "  -0.67%  "
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.224"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.868.29"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.266"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("E17").Value = "  +0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.14"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06722"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.007"
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.519.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.083.87"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.410"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1043"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +1.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.743"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.606"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02455"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06543"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.69%  "
$ws.Range("E37").Value = "  -0.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.818"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.025"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.186"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.63%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.243"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6368"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.008"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5989"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  +1.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.689"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("E49").Value = "  +0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.65%  "
$ws.Range("E51").Value = "  -10.05%  "
